# Daily attendance processing - 2025-11-06 23:21:34
#
# The "Recorded By" column (G) lists the users who recorded/edited a
# session, separated by ", ". This pass normalizes the ordering of
# those names by rotating the list so the last-listed recorder is
# moved to the front (most-recent-first ordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    "G2"   = "System, system, backup@backdoor.com"
    "G4"   = "System, backup@backdoor.com"
    "G5"   = "System, backup@backdoor.com"
    "G7"   = "System, admin@admin.com"
    "G8"   = "System, backup@backdoor.com"
    "G11"  = "System, dnasr281@gmail.com"
    "G17"  = "System, dnasr281@gmail.com"
    "G28"  = "System, system, backup@backdoor.com"
    "G30"  = "System, backup@backdoor.com"
    "G31"  = "System, backup@backdoor.com"
    "G33"  = "System, admin@admin.com"
    "G34"  = "System, backup@backdoor.com"
    "G37"  = "System, dnasr281@gmail.com"
    "G43"  = "System, dnasr281@gmail.com"
    "G54"  = "System, system, backup@backdoor.com"
    "G56"  = "System, backup@backdoor.com"
    "G57"  = "System, backup@backdoor.com"
    "G59"  = "System, admin@admin.com"
    "G60"  = "System, backup@backdoor.com"
    "G63"  = "System, dnasr281@gmail.com"
    "G69"  = "System, dnasr281@gmail.com"
    "G80"  = "System, backup@backdoor.com"
    "G81"  = "System, backup@backdoor.com"
    "G82"  = "System, backup@backdoor.com"
    "G87"  = "dnasr281@gmail.com, admin@admin.com"
    "G93"  = "System, dnasr281@gmail.com"
    "G94"  = "System, dnasr281@gmail.com"
    "G96"  = "System, dnasr281@gmail.com"
    "G106" = "System, backup@backdoor.com"
    "G107" = "System, backup@backdoor.com"
    "G108" = "System, backup@backdoor.com"
    "G113" = "dnasr281@gmail.com, admin@admin.com"
    "G119" = "System, dnasr281@gmail.com"
    "G120" = "System, dnasr281@gmail.com"
    "G122" = "System, dnasr281@gmail.com"
    "G132" = "System, backup@backdoor.com"
    "G133" = "System, backup@backdoor.com"
    "G134" = "System, backup@backdoor.com"
    "G139" = "dnasr281@gmail.com, admin@admin.com"
    "G145" = "System, dnasr281@gmail.com"
    "G146" = "System, dnasr281@gmail.com"
    "G148" = "System, dnasr281@gmail.com"
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
